$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, copying the header style from E1 (s="1")
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Fill in the time_taken values for each data row (F2:F53)
$ws.Range("F2").Value = "2021-10-05 13:41:40.309288"
$ws.Range("F3").Value = "2021-10-05 13:41:40.309300"
$ws.Range("F4").Value = "2021-10-05 13:41:40.309303"
$ws.Range("F5").Value = "2021-10-05 13:41:40.309306"
$ws.Range("F6").Value = "2021-10-05 13:41:40.309309"
$ws.Range("F7").Value = "2021-10-05 13:41:40.309312"
$ws.Range("F8").Value = "2021-10-05 13:41:40.309314"
$ws.Range("F9").Value = "2021-10-05 13:41:40.309317"
$ws.Range("F10").Value = "2021-10-05 13:41:40.309320"
$ws.Range("F11").Value = "2021-10-05 13:41:40.309322"
$ws.Range("F12").Value = "2021-10-05 13:41:40.309325"
$ws.Range("F13").Value = "2021-10-05 13:41:40.309328"
$ws.Range("F14").Value = "2021-10-05 13:41:40.309330"
$ws.Range("F15").Value = "2021-10-05 13:41:40.309333"
$ws.Range("F16").Value = "2021-10-05 13:41:40.309335"
$ws.Range("F17").Value = "2021-10-05 13:41:40.309338"
$ws.Range("F18").Value = "2021-10-05 13:41:40.309341"
$ws.Range("F19").Value = "2021-10-05 13:41:40.309343"
$ws.Range("F20").Value = "2021-10-05 13:41:40.309346"
$ws.Range("F21").Value = "2021-10-05 13:41:40.309349"
$ws.Range("F22").Value = "2021-10-05 13:41:40.309351"
$ws.Range("F23").Value = "2021-10-05 13:41:40.309354"
$ws.Range("F24").Value = "2021-10-05 13:41:40.309357"
$ws.Range("F25").Value = "2021-10-05 13:41:40.309359"
$ws.Range("F26").Value = "2021-10-05 13:41:40.309362"
$ws.Range("F27").Value = "2021-10-05 13:41:40.309365"
$ws.Range("F28").Value = "2021-10-05 13:41:40.309368"
$ws.Range("F29").Value = "2021-10-05 13:41:40.309370"
$ws.Range("F30").Value = "2021-10-05 13:41:40.309373"
$ws.Range("F31").Value = "2021-10-05 13:41:40.309376"
$ws.Range("F32").Value = "2021-10-05 13:41:40.309378"
$ws.Range("F33").Value = "2021-10-05 13:41:40.309381"
$ws.Range("F34").Value = "2021-10-05 13:41:40.309384"
$ws.Range("F35").Value = "2021-10-05 13:41:40.309387"
$ws.Range("F36").Value = "2021-10-05 13:41:40.309390"
$ws.Range("F37").Value = "2021-10-05 13:41:40.309392"
$ws.Range("F38").Value = "2021-10-05 13:41:40.309395"
$ws.Range("F39").Value = "2021-10-05 13:41:40.309398"
$ws.Range("F40").Value = "2021-10-05 13:41:40.309400"
$ws.Range("F41").Value = "2021-10-05 13:41:40.309403"
$ws.Range("F42").Value = "2021-10-05 13:41:40.309406"
$ws.Range("F43").Value = "2021-10-05 13:41:40.309409"
$ws.Range("F44").Value = "2021-10-05 13:41:40.309412"
$ws.Range("F45").Value = "2021-10-05 13:41:40.309414"
$ws.Range("F46").Value = "2021-10-05 13:41:40.309417"
$ws.Range("F47").Value = "2021-10-05 13:41:40.309419"
$ws.Range("F48").Value = "2021-10-05 13:41:40.309422"
$ws.Range("F49").Value = "2021-10-05 13:41:40.309425"
$ws.Range("F50").Value = "2021-10-05 13:41:40.309427"
$ws.Range("F51").Value = "2021-10-05 13:41:40.309430"
$ws.Range("F52").Value = "2021-10-05 13:41:40.309432"
$ws.Range("F53").Value = "2021-10-05 13:41:40.309435"
